# "Better cumulation of Kpis"
#
# Rename the first column header on the CapitalRemittance sheet from
# "Investor *" to "Stakeholder *", and move the active selection from
# K2:K6 up to A2 (matching the new topLeftCell/selection recorded in the
# saved sheet view).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CapitalRemittance")

# A1 holds the shared string "Investor *" -> rename to "Stakeholder *"
$ws.Range("A1").Value = "Stakeholder *"

# Move/collapse the selection to A2 (was K2:K6) and scroll the view back
# to the top-left of the sheet.
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
